# Update res_bus vm_pu values for the 380 kV case (rows 2-25, columns B-F,I-N)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.054076927776586
$ws.Range("D2").Value = 1.053932172296282
$ws.Range("E2").Value = 1.057759792043994
$ws.Range("F2").Value = 1.053056376451099
$ws.Range("I2").Value = 1.045355922328282
$ws.Range("J2").Value = 1.059090918719529
$ws.Range("K2").Value = 1.056676777510025
$ws.Range("L2").Value = 1.060493881222319
$ws.Range("M2").Value = 1.055803399478348
$ws.Range("N2").Value = 1.060594949303537

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.055721216404756
$ws.Range("D3").Value = 1.055234529281035
$ws.Range("E3").Value = 1.059355309705255
$ws.Range("F3").Value = 1.055324980377031
$ws.Range("I3").Value = 1.045922483563606
$ws.Range("J3").Value = 1.060382979144667
$ws.Range("K3").Value = 1.057791025940451
$ws.Range("L3").Value = 1.061901315339253
$ws.Range("M3").Value = 1.057881245810284
$ws.Range("N3").Value = 1.06188884460268

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.056781654681217
$ws.Range("D4").Value = 1.056074042083553
$ws.Range("E4").Value = 1.060384555887234
$ws.Range("F4").Value = 1.056788935499831
$ws.Range("I4").Value = 1.046285905230153
$ws.Range("J4").Value = 1.061215166299262
$ws.Range("K4").Value = 1.058508279537084
$ws.Range("L4").Value = 1.062808379309473
$ws.Range("M4").Value = 1.059221439483198
$ws.Range("N4").Value = 1.062722213558488

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.057226634739712
$ws.Range("D5").Value = 1.056426219310893
$ws.Range("E5").Value = 1.060816509248562
$ws.Range("F5").Value = 1.057403455257886
$ws.Range("I5").Value = 1.046437932491244
$ws.Range("J5").Value = 1.061564106348984
$ws.Range("K5").Value = 1.0588089289867
$ws.Range("L5").Value = 1.063188851466482
$ws.Range("M5").Value = 1.059783849243327
$ws.Range("N5").Value = 1.063071649143105

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.057301300746757
$ws.Range("D6").Value = 1.056485307549557
$ws.Range("E6").Value = 1.06088899305241
$ws.Range("F6").Value = 1.057506582239668
$ws.Range("I6").Value = 1.046463414436503
$ws.Range("J6").Value = 1.061622641853486
$ws.Range("K6").Value = 1.058859357842662
$ws.Range("L6").Value = 1.063252684530268
$ws.Range("M6").Value = 1.059878221937075
$ws.Range("N6").Value = 1.063130267774738

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.056787603758901
$ws.Range("D7").Value = 1.056078750839074
$ws.Range("E7").Value = 1.060390330563382
$ws.Range("F7").Value = 1.056797150349264
$ws.Range("I7").Value = 1.046287939584519
$ws.Range("J7").Value = 1.061219832416538
$ws.Range("K7").Value = 1.058512300287646
$ws.Range("L7").Value = 1.062813466547628
$ws.Range("M7").Value = 1.059228958356051
$ws.Range("N7").Value = 1.062726886302185

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.054633362687224
$ws.Range("D8").Value = 1.054372979596288
$ws.Range("E8").Value = 1.058299668099346
$ws.Range("F8").Value = 1.053823902655892
$ws.Range("I8").Value = 1.045548056947521
$ws.Range("J8").Value = 1.059528384312246
$ws.Range("K8").Value = 1.05705412415025
$ws.Range("L8").Value = 1.060970292265811
$ws.Range("M8").Value = 1.056506525055927
$ws.Range("N8").Value = 1.061033036147574

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.050809582285399
$ws.Range("D9").Value = 1.051342150049948
$ws.Range("E9").Value = 1.054590776442883
$ws.Range("F9").Value = 1.048552950293224
$ws.Range("I9").Value = 1.04421963406032
$ws.Range("J9").Value = 1.056517658361932
$ws.Range("K9").Value = 1.054455478806698
$ws.Range("L9").Value = 1.057693877952693
$ws.Range("M9").Value = 1.051675115240479
$ws.Range("N9").Value = 1.058018034621075

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.048240730545649
$ws.Range("D10").Value = 1.049303998569707
$ws.Range("E10").Value = 1.052100509819019
$ws.Range("F10").Value = 1.0450159716586
$ws.Range("I10").Value = 1.043317033289626
$ws.Range("J10").Value = 1.054489376365819
$ws.Range("K10").Value = 1.052702727374238
$ws.Range("L10").Value = 1.055489552534981
$ws.Range("M10").Value = 1.048429663514031
$ws.Range("N10").Value = 1.05598687223188

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.047123496545224
$ws.Range("D11").Value = 1.048417111306795
$ws.Range("E11").Value = 1.051017792934777
$ws.Range("F11").Value = 1.043478584487077
$ws.Range("I11").Value = 1.042922076483735
$ws.Range("J11").Value = 1.053605906909215
$ws.Range("K11").Value = 1.051938784817433
$ws.Range("L11").Value = 1.054530101528351
$ws.Range("M11").Value = 1.047018183311854
$ws.Range("N11").Value = 1.055102148147313

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.046707749446566
$ws.Range("D12").Value = 1.048087012838483
$ws.Range("E12").Value = 1.050614941936701
$ws.Range("F12").Value = 1.042906618462265
$ws.Range("I12").Value = 1.042774744223619
$ws.Range("J12").Value = 1.053276947861128
$ws.Range("K12").Value = 1.051654258814913
$ws.Range("L12").Value = 1.05417295575371
$ws.Range("M12").Value = 1.046492938218847
$ws.Range("N12").Value = 1.054772721939646

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.046796963171662
$ws.Range("D13").Value = 1.048157850597447
$ws.Range("E13").Value = 1.050701385992051
$ws.Range("F13").Value = 1.043029348915945
$ws.Range("I13").Value = 1.042806376029075
$ws.Range("J13").Value = 1.053347547123028
$ws.Range("K13").Value = 1.051715325396587
$ws.Range("L13").Value = 1.054249599539041
$ws.Range("M13").Value = 1.046605648939257
$ws.Range("N13").Value = 1.054843421460595

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.047089146311427
$ws.Range("D14").Value = 1.048389838999656
$ws.Range("E14").Value = 1.050984507160422
$ws.Range("F14").Value = 1.043431324392303
$ws.Range("I14").Value = 1.042909910805214
$ws.Range("J14").Value = 1.053578731449988
$ws.Range("K14").Value = 1.051915281471212
$ws.Range("L14").Value = 1.054500595411719
$ws.Range("M14").Value = 1.046974786118504
$ws.Range("N14").Value = 1.055074934095816

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.047269069126077
$ws.Range("D15").Value = 1.048532685606572
$ws.Range("E15").Value = 1.051158856563559
$ws.Range("F15").Value = 1.043678872916177
$ws.Range("I15").Value = 1.042973618577285
$ws.Range("J15").Value = 1.053721065353026
$ws.Range("K15").Value = 1.052038379382722
$ws.Range("L15").Value = 1.054655140634744
$ws.Range("M15").Value = 1.047202095793558
$ws.Range("N15").Value = 1.055217470129323

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.04831477267264
$ws.Range("D16").Value = 1.049362765352801
$ws.Range("E16").Value = 1.0521722715859
$ws.Range("F16").Value = 1.045117876898703
$ws.Range("I16").Value = 1.043343157686608
$ws.Range("J16").Value = 1.05454789815985
$ws.Range("K16").Value = 1.052753321402098
$ws.Range("L16").Value = 1.055553122096422
$ws.Range("M16").Value = 1.048523206062693
$ws.Range("N16").Value = 1.056045477133572

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.048969387952311
$ws.Range("D17").Value = 1.04988227625955
$ws.Range("E17").Value = 1.052806765097235
$ws.Range("F17").Value = 1.04601893632006
$ws.Range("I17").Value = 1.043573849872655
$ws.Range("J17").Value = 1.05506514215935
$ws.Range("K17").Value = 1.053200439894283
$ws.Range("L17").Value = 1.056115061165875
$ws.Range("M17").Value = 1.049350228298089
$ws.Range("N17").Value = 1.05656345567888

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.049350742134574
$ws.Range("D18").Value = 1.050184879078251
$ws.Range("E18").Value = 1.053176429701733
$ws.Range("F18").Value = 1.046543946630136
$ws.Range("I18").Value = 1.043708011294345
$ws.Range("J18").Value = 1.055366340267148
$ws.Range("K18").Value = 1.053460756081725
$ws.Range("L18").Value = 1.056442352982767
$ws.Range("M18").Value = 1.049832022337206
$ws.Range("N18").Value = 1.056865081522539

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.04948069462049
$ws.Range("D19").Value = 1.050287988179176
$ws.Range("E19").Value = 1.053302404292007
$ws.Range("F19").Value = 1.046722867169568
$ws.Range("I19").Value = 1.04375368970259
$ws.Range("J19").Value = 1.05546895645726
$ws.Range("K19").Value = 1.05354943610266
$ws.Range("L19").Value = 1.056553870520059
$ws.Range("M19").Value = 1.049996201805246
$ws.Range("N19").Value = 1.056967843439411

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.04889920287514
$ws.Range("D20").Value = 1.049826581112977
$ws.Range("E20").Value = 1.052738734012533
$ws.Range("F20").Value = 1.045922319579929
$ws.Range("I20").Value = 1.043549139949976
$ws.Range("J20").Value = 1.055009698771003
$ws.Range("K20").Value = 1.053152518109794
$ws.Range("L20").Value = 1.056054819960454
$ws.Range("M20").Value = 1.049261558301498
$ws.Range("N20").Value = 1.056507933554562

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.047003126707017
$ws.Range("D21").Value = 1.048321542786927
$ws.Range("E21").Value = 1.050901154001578
$ws.Range("F21").Value = 1.043312978073938
$ws.Range("I21").Value = 1.042879439787239
$ws.Range("J21").Value = 1.053510675621373
$ws.Range("K21").Value = 1.051856420587592
$ws.Range("L21").Value = 1.054426704583876
$ws.Range("M21").Value = 1.046866111124398
$ws.Range("N21").Value = 1.055006781620119

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.045806600623805
$ws.Range("D22").Value = 1.047371387018978
$ws.Range("E22").Value = 1.049741840868614
$ws.Range("F22").Value = 1.041667088522514
$ws.Range("I22").Value = 1.042454736881605
$ws.Range("J22").Value = 1.0525635483087
$ws.Range("K22").Value = 1.051037086850984
$ws.Range("L22").Value = 1.053398620918956
$ws.Range("M22").Value = 1.04535443751114
$ws.Range("N22").Value = 1.054058309278025

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.046441324339207
$ws.Range("D23").Value = 1.047875455220655
$ws.Range("E23").Value = 1.050356795674966
$ws.Range("F23").Value = 1.042540118487592
$ws.Range("I23").Value = 1.04268022727064
$ws.Range("J23").Value = 1.053066083013248
$ws.Range("K23").Value = 1.051471855434215
$ws.Range("L23").Value = 1.053944052576882
$ws.Range("M23").Value = 1.046156341629046
$ws.Range("N23").Value = 1.054561557639495

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.048930917947962
$ws.Range("D24").Value = 1.049851748647103
$ws.Range("E24").Value = 1.052769475643373
$ws.Range("F24").Value = 1.045965978262756
$ws.Range("I24").Value = 1.043560306528288
$ws.Range("J24").Value = 1.055034752801133
$ws.Range("K24").Value = 1.053174173385772
$ws.Range("L24").Value = 1.056082041843958
$ws.Range("M24").Value = 1.04930162629082
$ws.Range("N24").Value = 1.056533023164288

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.051801511731259
$ws.Range("D25").Value = 1.052128737565832
$ws.Range("E25").Value = 1.055552658580613
$ws.Range("F25").Value = 1.049919548310517
$ws.Range("I25").Value = 1.044566026727298
$ws.Range("J25").Value = 1.057299665394662
$ws.Range("K25").Value = 1.055130819258451
$ws.Range("L25").Value = 1.058544377923501
$ws.Range("M25").Value = 1.052928355143053
$ws.Range("N25").Value = 1.058801152193489
